$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U71"), 0, 0)
